$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - 11.03.2020, 10:40 -> 11:35
$ws.Range("C8").Value = 43901
$ws.Range("D8").Value = 0.44444444444444442
$ws.Range("E8").Value = 0.4826388888888889
$ws.Range("G8").Value = "Bataille Navale/ICT 431"
$ws.Range("H8").Value = "Théorie"
$ws.Range("I8").Value = "Théorie sur la planification du projet"
$ws.Range("J8").Value = "GitHub"

# Row 9 - 11.03.2020, 11:40 -> 12:15
$ws.Range("C9").Value = 43901
$ws.Range("D9").Value = 0.4861111111111111
$ws.Range("E9").Value = 0.51041666666666663
$ws.Range("G9").Value = "Bataille Navale/ICT 431"
$ws.Range("H9").Value = "Planification"
$ws.Range("I9").Value = "Realisation des sprints pour planifier le projet"
$ws.Range("J9").Value = "GitHub"

# Update selection to match the author's last active cell
$ws.Range("K9").Select()
